# Sidste opdatering af tidsregistrering
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- Append two new time-registration rows for 23-3-2017 ---
# Clone number formats from the last existing data row (row 42) so the new
# rows reuse the same date/time styles instead of minting new ones.
$ws.Range("A42").Copy()
$ws.Range("A44:A45").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G42:H42").Copy()
$ws.Range("G44:H45").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 44
$ws.Range("A44").Value = "3/23/2017"
$ws.Range("B44").Value = "Tolga Yasin Kücük"
$ws.Range("D44").Value = "Test Designer"
$ws.Range("F44").Value = "Test for flydespaending"
$ws.Range("G44").Value = 0.34375
$ws.Range("H44").Value = 0.3923611111111111

# --- Fix role for three existing rows (Test Designer -> Test Analysist) ---
$ws.Range("D25").Value = "Test Analysist"
$ws.Range("D26").Value = "Test Analysist"
$ws.Range("D33").Value = "Test Analysist"

# Row 45
$ws.Range("A45").Value = "3/23/2017"
$ws.Range("B45").Value = "Tolga Yasin Kücük"
$ws.Range("D45").Value = "Test Designer"
$ws.Range("F45").Value = "Test for dimensionerende kraft(Ikke færdig)"
$ws.Range("G45").Value = 0.39583333333333331
$ws.Range("H45").Value = 0.43402777777777773
$ws.Range("I45").Value = 2

# --- Fix role for three existing rows ("?" -> User-Interface Designer) ---
$ws.Range("D40").Value = "User-Interface Designer"
$ws.Range("D41").Value = "User-Interface Designer"
$ws.Range("D42").Value = "User-Interface Designer"

# --- Update the view: scroll position + active selection ---
$ws.Activate()
$ws.Range("I43").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 34
